# Update need_to_buy.xlsx data: shift one day forward and refresh values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 46010
$ws.Range("B2").Value = 10536.1356529996
$ws.Range("C2").Value = 8893.80146086179
$ws.Range("D2").Value = 15792.26
$ws.Range("E2").Value = 6602.34921778274
$ws.Range("F2").Value = -12.337888389811

# Row 3
$ws.Range("A3").Value = 46011
$ws.Range("B3").Value = 4140.55085362935
$ws.Range("C3").Value = 6277.90749539817
$ws.Range("D3").Value = 11232.26
$ws.Range("E3").Value = 6743.65016178358
$ws.Range("F3").Value = 74.5540690492393

# Row 4
$ws.Range("A4").Value = 46012
$ws.Range("B4").Value = 4123.91555761858
$ws.Range("C4").Value = 6435.05794616579
$ws.Range("D4").Value = 11232.26
$ws.Range("E4").Value = 6880.94764355747
$ws.Range("F4").Value = 86.8227329051358

# Row 5
$ws.Range("A5").Value = 46013
$ws.Range("B5").Value = 9917.52662600306
$ws.Range("C5").Value = 9324.53582888978
$ws.Range("D5").Value = 11232.26
$ws.Range("E5").Value = 7386.07672865471
$ws.Range("F5").Value = 228.264689897687

# Row 6
$ws.Range("A6").Value = 46014
$ws.Range("B6").Value = 9919.8033347958
$ws.Range("C6").Value = 9412.16682424713
$ws.Range("D6").Value = 17376.26
$ws.Range("E6").Value = 7417.08993717028
$ws.Range("F6").Value = -22.7918016076078

# Row 7
$ws.Range("A7").Value = 46015
$ws.Range("B7").Value = 9566.18303866457
$ws.Range("C7").Value = 9325.34374210863
$ws.Range("D7").Value = 17208.26
$ws.Range("E7").Value = 8113.18612739547
$ws.Range("F7").Value = 9.59457789600432

# Row 8
$ws.Range("A8").Value = 46016
$ws.Range("B8").Value = 8471.80799626025
$ws.Range("C8").Value = 8692.32427628313
$ws.Range("D8").Value = 16416.26
$ws.Range("E8").Value = 7768.80965290935
$ws.Range("F8").Value = 1.86974704968679

# Row 9
$ws.Range("A9").Value = 46017
$ws.Range("B9").Value = 8471.80799626025
$ws.Range("C9").Value = 8354.68573219736
$ws.Range("D9").Value = 15840.26
$ws.Range("E9").Value = 7768.80965290935
$ws.Range("F9").Value = 11.8014743794461

# Row 10
$ws.Range("A10").Value = 46018
$ws.Range("B10").Value = 8547.24893097507
$ws.Range("C10").Value = 9421.06739776354
$ws.Range("D10").Value = 16632.26
$ws.Range("E10").Value = 7775.98859096
$ws.Range("F10").Value = 23.5331661968142

# Row 11
$ws.Range("A11").Value = 46019
$ws.Range("B11").Value = 8471.80799626025
$ws.Range("C11").Value = 9617.5757544937
$ws.Range("D11").Value = 17448.26
$ws.Range("E11").Value = 7768.80965290935
$ws.Range("F11").Value = -2.57810802487287

# Row 12
$ws.Range("A12").Value = 46020
$ws.Range("B12").Value = 9603.52690270833
$ws.Range("C12").Value = 10675.5261421174
$ws.Range("D12").Value = 19176.26
$ws.Range("E12").Value = 8403.11249713734
$ws.Range("F12").Value = -4.0675566977205

# Row 13
$ws.Range("A13").Value = 46021
$ws.Range("B13").Value = 9603.52690270833
$ws.Range("C13").Value = 11413.6138794642
$ws.Range("D13").Value = 11232.26
$ws.Range("E13").Value = 8403.11249713734
$ws.Range("F13").Value = 357.686099025066

# Row 14
$ws.Range("A14").Value = 46022
$ws.Range("B14").Value = 9603.52690270833
$ws.Range("C14").Value = 11509.9711885105
$ws.Range("D14").Value = 11232.26
$ws.Range("E14").Value = 8403.11249713734
$ws.Range("F14").Value = 361.700986901992

# Row 15
$ws.Range("A15").Value = 46023
$ws.Range("B15").Value = 4937.72735627544
$ws.Range("C15").Value = 9523.56449913079
$ws.Range("D15").Value = 10875.86
$ws.Range("E15").Value = 8491.08806734464
$ws.Range("F15").Value = 297.44969026981
